$wb = $excel.ActiveWorkbook

# Neodymium
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C2").Value = [double]"0.0002735394444618311"
$ws.Range("D2").Value = [double]"0.728535354108951"
$ws.Range("E2").Value = [double]"0.8339924611837429"
$ws.Range("B3").Value = [double]"3.676322720692059E-10"
$ws.Range("C3").Value = [double]"0.0132372580562949"
$ws.Range("D3").Value = [double]"0.6352624942659592"
$ws.Range("E3").Value = [double]"0.7415669097852163"
$ws.Range("B4").Value = [double]"5.738846834589291E-12"
$ws.Range("C4").Value = [double]"0.01197210825341109"
$ws.Range("D4").Value = [double]"0.5191851368327279"
$ws.Range("E4").Value = [double]"0.6548198374626368"
$ws.Range("C5").Value = [double]"2.664856196090391E-07"
$ws.Range("D5").Value = [double]"0.028701895446782"
$ws.Range("E5").Value = [double]"0.05468359805677819"

# Dysprosium
$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C2").Value = [double]"0.0003224889070000346"
$ws.Range("D2").Value = [double]"0.8589056343216422"
$ws.Range("E2").Value = [double]"0.9832341283815936"
$ws.Range("B3").Value = [double]"4.334195012744524E-10"
$ws.Range("C3").Value = [double]"0.0156060450098911"
$ws.Range("D3").Value = [double]"0.748941739780901"
$ws.Range("E3").Value = [double]"0.8742691668272271"
$ws.Range("B4").Value = [double]"6.765804642063433E-12"
$ws.Range("C4").Value = [double]"0.01411449859717533"
$ws.Range("D4").Value = [double]"0.6120925179081902"
$ws.Range("E4").Value = [double]"0.7719988394387939"
$ws.Range("C5").Value = [double]"3.141728110475074E-07"
$ws.Range("D5").Value = [double]"0.03383805545731329"
$ws.Range("E5").Value = [double]"0.06446914375678113"

# Copper
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"1.254977486862241E-05"
$ws.Range("C2").Value = [double]"0.009674324557565138"
$ws.Range("D2").Value = [double]"1.396182838935287"
$ws.Range("E2").Value = [double]"1.274673889735591"
$ws.Range("B3").Value = [double]"8.533411143136702E-05"
$ws.Range("C3").Value = [double]"0.03490138443766061"
$ws.Range("D3").Value = [double]"0.9914554302842176"
$ws.Range("E3").Value = [double]"0.9787974795986573"
$ws.Range("B4").Value = [double]"0.0002530745263323527"
$ws.Range("C4").Value = [double]"0.009342366228943881"
$ws.Range("D4").Value = [double]"0.8414105811741163"
$ws.Range("E4").Value = [double]"0.9858045873434578"
$ws.Range("B5").Value = [double]"7.95025142153212E-05"
$ws.Range("C5").Value = [double]"0.02048333445463943"
$ws.Range("D5").Value = [double]"1.208410145871788"
$ws.Range("E5").Value = [double]"0.995634747265968"

# Raw silicon
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"0.000107623096286272"
$ws.Range("C2").Value = [double]"0.007486642629472633"
$ws.Range("D2").Value = [double]"2.016567546390563"
$ws.Range("E2").Value = [double]"1.868247443550294"
$ws.Range("B3").Value = [double]"0.0001148540841621532"
$ws.Range("C3").Value = [double]"0.02501457324738229"
$ws.Range("D3").Value = [double]"1.062027883849227"
$ws.Range("E3").Value = [double]"1.038192857648591"
$ws.Range("B4").Value = [double]"0.0007359536893654886"
$ws.Range("C4").Value = [double]"0.007021125419683741"
$ws.Range("D4").Value = [double]"1.097468629610411"
$ws.Range("E4").Value = [double]"1.303346366569932"
$ws.Range("B5").Value = [double]"0.0003951808733664171"
$ws.Range("C5").Value = [double]"0.008916107934033045"
$ws.Range("D5").Value = [double]"1.879760438469067"
$ws.Range("E5").Value = [double]"1.549935507400286"

